# Auto-generated script to apply scheduled market-data refresh updates
# to the leve profit tables across all job sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value2 = 8336.385
$ws.Range("I6").Value2 = 11326.777
$ws.Range("K6").Value2 = 33980.331
$ws.Range("M6").Value2 = -33868.331
$ws.Range("H33").Value2 = 408.3846
$ws.Range("I33").Value2 = 309.16666
$ws.Range("K33").Value2 = 309.16666
$ws.Range("M33").Value2 = -80.16665999999998
$ws.Range("H69").Value2 = 0
$ws.Range("I69").Value2 = 0
$ws.Range("J69").Value2 = 0
$ws.Range("K69").Value2 = 0
$ws.Range("L69").ClearContents()
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value2 = 0
$ws.Range("H72").Value2 = 0
$ws.Range("I72").Value2 = 0
$ws.Range("J72").Value2 = 0
$ws.Range("K72").Value2 = 0
$ws.Range("L72").ClearContents()
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value2 = 0
$ws.Range("H76").Value2 = 6248.25
$ws.Range("I76").Value2 = 4997.6665
$ws.Range("K76").Value2 = 4997.6665
$ws.Range("M76").Value2 = -4682.6665
$ws.Range("H79").Value2 = 6248.25
$ws.Range("I79").Value2 = 4997.6665
$ws.Range("K79").Value2 = 4997.6665
$ws.Range("M79").Value2 = -3905.6665
$ws.Range("H86").Value2 = 62500520
$ws.Range("I86").Value2 = 90909430
$ws.Range("J86").Value2 = 903
$ws.Range("K86").Value2 = 90909430
$ws.Range("L86").Value2 = 903
$ws.Range("M86").Value2 = -90908307
$ws.Range("N86").Value2 = -3149
$ws.Range("H89").Value2 = 62500520
$ws.Range("I89").Value2 = 90909430
$ws.Range("J89").Value2 = 903
$ws.Range("K89").Value2 = 454547150
$ws.Range("L89").Value2 = 4515
$ws.Range("M89").Value2 = -454541534
$ws.Range("N89").Value2 = -15747
$ws.Range("H98").Value2 = 420.3243
$ws.Range("I98").Value2 = 453.52942
$ws.Range("K98").Value2 = 453.52942
$ws.Range("M98").Value2 = 1044.47058
$ws.Range("H106").Value2 = 1493.4762
$ws.Range("I106").Value2 = 1397.9445
$ws.Range("K106").Value2 = 1397.9445
$ws.Range("M106").Value2 = -766.9445000000001
$ws.Range("H122").Value2 = 420.3243
$ws.Range("I122").Value2 = 453.52942
$ws.Range("K122").Value2 = 1360.58826
$ws.Range("M122").Value2 = 1089.41174
$ws.Range("H132").Value2 = 120314.23
$ws.Range("I132").Value2 = 154634
$ws.Range("K132").Value2 = 463902
$ws.Range("M132").Value2 = -461372
$ws.Range("H138").Value2 = 5528.8677
$ws.Range("I138").Value2 = 1697.619
$ws.Range("J138").Value2 = 6826.5483
$ws.Range("K138").Value2 = 5092.857
$ws.Range("L138").Value2 = 20479.6449
$ws.Range("M138").Value2 = 47.14300000000003
$ws.Range("N138").Value2 = -30759.6449
$ws.Range("H141").Value2 = 3123.18
$ws.Range("I141").Value2 = 2971.3333
$ws.Range("K141").Value2 = 8913.999899999999
$ws.Range("M141").Value2 = -3733.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 216472.81
$ws.Range("I61").Value2 = 2253.0881
$ws.Range("K61").Value2 = 2253.0881
$ws.Range("M61").Value2 = -2041.0881
$ws.Range("H74").Value2 = 11820403
$ws.Range("I74").Value2 = 20834488
$ws.Range("J74").Value2 = 1003501.8
$ws.Range("K74").Value2 = 20834488
$ws.Range("L74").Value2 = 1003501.8
$ws.Range("M74").Value2 = -20833614
$ws.Range("N74").Value2 = -1005249.8
$ws.Range("H77").Value2 = 11820403
$ws.Range("I77").Value2 = 20834488
$ws.Range("J77").Value2 = 1003501.8
$ws.Range("K77").Value2 = 104172440
$ws.Range("L77").Value2 = 5017509
$ws.Range("M77").Value2 = -104168072
$ws.Range("N77").Value2 = -5026245
$ws.Range("H130").Value2 = 84924.5
$ws.Range("J130").Value2 = 84924.5
$ws.Range("L130").Value2 = 84924.5
$ws.Range("N130").Value2 = -94964.5
$ws.Range("H132").Value2 = 17699.36
$ws.Range("I132").Value2 = 23387.885
$ws.Range("K132").Value2 = 70163.655
$ws.Range("M132").Value2 = -67633.655
$ws.Range("H136").Value2 = 216472.81
$ws.Range("I136").Value2 = 2253.0881
$ws.Range("K136").Value2 = 6759.2643
$ws.Range("M136").Value2 = -4209.2643

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 2610.7693
$ws.Range("I20").Value2 = 2049.1428
$ws.Range("K20").Value2 = 2049.1428
$ws.Range("M20").Value2 = -1802.1428
$ws.Range("H94").Value2 = 1477.8422
$ws.Range("J94").Value2 = 1542.375
$ws.Range("L94").Value2 = 1542.375
$ws.Range("N94").Value2 = -2444.375
$ws.Range("H134").Value2 = 3038.5667
$ws.Range("I134").Value2 = 2009
$ws.Range("K134").Value2 = 6027
$ws.Range("M134").Value2 = -3492

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 37041876
$ws.Range("I31").Value2 = 90910760
$ws.Range("J31").Value2 = 7016.3125
$ws.Range("K31").Value2 = 90910760
$ws.Range("L31").Value2 = 7016.3125
$ws.Range("M31").Value2 = -90910465
$ws.Range("N31").Value2 = -7606.3125
$ws.Range("H34").Value2 = 37041876
$ws.Range("I34").Value2 = 90910760
$ws.Range("J34").Value2 = 7016.3125
$ws.Range("K34").Value2 = 90910760
$ws.Range("L34").Value2 = 7016.3125
$ws.Range("M34").Value2 = -90910558
$ws.Range("N34").Value2 = -7420.3125
$ws.Range("H58").Value2 = 287481.62
$ws.Range("I58").Value2 = 1547.5385
$ws.Range("K58").Value2 = 1547.5385
$ws.Range("M58").Value2 = -1344.5385
$ws.Range("H60").Value2 = 14801.143
$ws.Range("J60").Value2 = 16752.5
$ws.Range("L60").Value2 = 16752.5
$ws.Range("N60").Value2 = -17774.5
$ws.Range("H86").Value2 = 6251.2
$ws.Range("J86").Value2 = 6683.75
$ws.Range("L86").Value2 = 6683.75
$ws.Range("N86").Value2 = -8929.75
$ws.Range("H89").Value2 = 6251.2
$ws.Range("J89").Value2 = 6683.75
$ws.Range("L89").Value2 = 33418.75
$ws.Range("N89").Value2 = -44650.75
$ws.Range("H134").Value2 = 1954.3636
$ws.Range("I134").Value2 = 1768.6571
$ws.Range("J134").Value2 = 2676.5557
$ws.Range("K134").Value2 = 5305.971299999999
$ws.Range("L134").Value2 = 8029.6671
$ws.Range("M134").Value2 = -2770.971299999999
$ws.Range("N134").Value2 = -13099.6671
$ws.Range("H136").Value2 = 287481.62
$ws.Range("I136").Value2 = 1547.5385
$ws.Range("K136").Value2 = 4642.6155
$ws.Range("M136").Value2 = -2092.6155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value2 = 1656.9333
$ws.Range("I68").Value2 = 1099.6666
$ws.Range("J68").Value2 = 1796.25
$ws.Range("K68").Value2 = 3298.9998
$ws.Range("L68").Value2 = 5388.75
$ws.Range("M68").Value2 = -2487.9998
$ws.Range("N68").Value2 = -7010.75
$ws.Range("H71").Value2 = 1656.9333
$ws.Range("I71").Value2 = 1099.6666
$ws.Range("J71").Value2 = 1796.25
$ws.Range("K71").Value2 = 9896.999400000001
$ws.Range("L71").Value2 = 16166.25
$ws.Range("M71").Value2 = -5840.999400000001
$ws.Range("N71").Value2 = -24278.25
$ws.Range("H131").Value2 = 16004748
$ws.Range("I131").Value2 = 16767478
$ws.Range("J131").Value2 = 15722255
$ws.Range("K131").Value2 = 50302434
$ws.Range("L131").Value2 = 47166765
$ws.Range("M131").Value2 = -50297394
$ws.Range("N131").Value2 = -47176845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value2 = 1488.1305
$ws.Range("I97").Value2 = 1248.6875
$ws.Range("J97").Value2 = 2035.4286
$ws.Range("K97").Value2 = 1248.6875
$ws.Range("L97").Value2 = 2035.4286
$ws.Range("M97").Value2 = -752.6875
$ws.Range("N97").Value2 = -3027.4286
$ws.Range("H122").Value2 = 260284.53
$ws.Range("I122").Value2 = 358022.47
$ws.Range("J122").Value2 = 7794.8335
$ws.Range("K122").Value2 = 1074067.41
$ws.Range("L122").Value2 = 23384.5005
$ws.Range("M122").Value2 = -1071617.41
$ws.Range("N122").Value2 = -28284.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 3874.2898
$ws.Range("I7").Value2 = 3151.8298
$ws.Range("K7").Value2 = 3151.8298
$ws.Range("M7").Value2 = -3039.8298
$ws.Range("H16").Value2 = 66715.5
$ws.Range("I16").Value2 = 102394.6
$ws.Range("J16").Value2 = 7250.3335
$ws.Range("K16").Value2 = 102394.6
$ws.Range("L16").Value2 = 7250.3335
$ws.Range("M16").Value2 = -102224.6
$ws.Range("N16").Value2 = -7590.3335
$ws.Range("H40").Value2 = 21829032
$ws.Range("I40").Value2 = 8335790.5
$ws.Range("J40").Value2 = 55562136
$ws.Range("K40").Value2 = 8335790.5
$ws.Range("L40").Value2 = 55562136
$ws.Range("M40").Value2 = -8335654.5
$ws.Range("N40").Value2 = -55562408
$ws.Range("H93").Value2 = 532.1667
$ws.Range("I93").Value2 = 448.25
$ws.Range("K93").Value2 = 448.25
$ws.Range("M93").Value2 = 799.75
$ws.Range("H95").Value2 = 49999
$ws.Range("J95").Value2 = 49999
$ws.Range("L95").Value2 = 49999
$ws.Range("N95").Value2 = -55491
$ws.Range("H100").Value2 = 2298.3333
$ws.Range("I100").Value2 = 2298.125
$ws.Range("K100").Value2 = 2298.125
$ws.Range("M100").Value2 = -1757.125
$ws.Range("H126").Value2 = 3874.2898
$ws.Range("I126").Value2 = 3151.8298
$ws.Range("K126").Value2 = 9455.4894
$ws.Range("M126").Value2 = -6985.4894
$ws.Range("H136").Value2 = 3057.5312
$ws.Range("I136").Value2 = 2145.6
$ws.Range("K136").Value2 = 6436.799999999999
$ws.Range("M136").Value2 = -3886.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value2 = 2415.5
$ws.Range("I96").Value2 = 869
$ws.Range("K96").Value2 = 869
$ws.Range("M96").Value2 = 504
$ws.Range("H122").Value2 = 5352.077
$ws.Range("I122").Value2 = 3371.8572
$ws.Range("J122").Value2 = 7662.3335
$ws.Range("K122").Value2 = 10115.5716
$ws.Range("L122").Value2 = 22987.0005
$ws.Range("M122").Value2 = -7665.571599999999
$ws.Range("N122").Value2 = -27887.0005
$ws.Range("H126").Value2 = 5761.231
$ws.Range("I126").Value2 = 5872.364
$ws.Range("K126").Value2 = 17617.092
$ws.Range("M126").Value2 = -15147.092
$ws.Range("H132").Value2 = 315720.72
$ws.Range("J132").Value2 = 1440679
$ws.Range("L132").Value2 = 4322037
$ws.Range("N132").Value2 = -4327097
$ws.Range("H136").Value2 = 115037.25
$ws.Range("J136").Value2 = 189984.11
$ws.Range("L136").Value2 = 569952.33
$ws.Range("N136").Value2 = -575052.33

